# Append a new customer record (phone 71076783) with 0 points.
# Mirrors the target diff: a new row 10 is added to sheetData with
#   A10 = "71076783" (stored as text, like the diff's inlineStr cell)
#   B10 = "" (empty, birthday unknown)
#   C10 = 0  (numeric total_points)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10

$phoneCell = $ws.Cells.Item($lastRow, 1)
# Force text storage so the long phone number isn't coerced into a
# numeric cell (the source data keeps this value as a string).
$phoneCell.NumberFormat = "@"
$phoneCell.Value = "71076783"
$phoneCell.ClearFormats()

$birthdayCell = $ws.Cells.Item($lastRow, 2)
$birthdayCell.NumberFormat = "@"
$birthdayCell.Value = ""
$birthdayCell.ClearFormats()

$pointsCell = $ws.Cells.Item($lastRow, 3)
$pointsCell.Value = 0
